$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

# Row 18 previously held the "juniorAssistant" approval entry; replace it
# with the "LightingSuperintendent" / "Lighting Superintendent" entry for
# JAGADEESH MADARAPU, as part of completing Forward/Close Grievance.
$ws.Range("A18").Value = "LightingSuperintendent"
$ws.Range("B18").Value = "ENGINEERING"
$ws.Range("C18").Value = "Lighting Superintendent"
$ws.Range("D18").Value = "JAGADEESH MADARAPU"
$ws.Range("E18").Value = "Sanctioned and shall grievance be processed"

# Update the active selection left behind in the saved sheet view.
$ws.Range("B21").Select()
